# Update TPM-derived NATMI ligand-receptor metrics for Slit3-Robo4 (OldD0)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.232451333333333
$ws.Range("H2").Value = 6.697354
$ws.Range("I2").Value = 0.01414074962829973
$ws.Range("J2").Value = 0.01414074962829973
$ws.Range("M2").Value = 39.327127
$ws.Range("N2").Value = 117.981381
$ws.Range("O2").Value = 0.9923865713449503
$ws.Range("P2").Value = 0.9923865713449502
$ws.Range("Q2").Value = 87.79589710731932
$ws.Range("R2").Value = 790.163073965874
$ws.Range("S2").Value = 0.01403309003987575
$ws.Range("T2").Value = 0.01403309003987575

$ws.Range("G3").Value = 2.232451333333333
$ws.Range("H3").Value = 6.697354
$ws.Range("I3").Value = 0.01414074962829973
$ws.Range("J3").Value = 0.01414074962829973
$ws.Range("O3").Value = 0.001455135597170125
$ws.Range("P3").Value = 0.001455135597170125
$ws.Range("Q3").Value = 0.1287350502871111
$ws.Range("R3").Value = 1.158615452584
$ws.Range("S3").Value = 0.00002057670815480916
$ws.Range("T3").Value = 0.00002057670815480915

$ws.Range("G4").Value = 2.232451333333333
$ws.Range("H4").Value = 6.697354
$ws.Range("I4").Value = 0.01414074962829973
$ws.Range("J4").Value = 0.01414074962829973
$ws.Range("O4").Value = 0.00615829305787961
$ws.Range("P4").Value = 0.006158293057879609
$ws.Range("Q4").Value = 0.5448208180946665
$ws.Range("R4").Value = 4.903387362851999
$ws.Range("S4").Value = 0.00008708288026917192
$ws.Range("T4").Value = 0.0000870828802691719

$ws.Range("I5").Value = 0.8099327614075106
$ws.Range("J5").Value = 0.8099327614075106
$ws.Range("M5").Value = 39.327127
$ws.Range("N5").Value = 117.981381
$ws.Range("O5").Value = 0.9923865713449503
$ws.Range("P5").Value = 0.9923865713449502
$ws.Range("Q5").Value = 5028.642416670158
$ws.Range("R5").Value = 45257.78175003143
$ws.Range("S5").Value = 0.8037663961131472
$ws.Range("T5").Value = 0.803766396113147

$ws.Range("I6").Value = 0.8099327614075106
$ws.Range("J6").Value = 0.8099327614075106
$ws.Range("O6").Value = 0.001455135597170125
$ws.Range("P6").Value = 0.001455135597170125
$ws.Range("S6").Value = 0.001178561992438366
$ws.Range("T6").Value = 0.001178561992438366

$ws.Range("I7").Value = 0.8099327614075106
$ws.Range("J7").Value = 0.8099327614075106
$ws.Range("O7").Value = 0.00615829305787961
$ws.Range("P7").Value = 0.006158293057879609
$ws.Range("S7").Value = 0.004987803301925135
$ws.Range("T7").Value = 0.004987803301925134

$ws.Range("G8").Value = 27.77415166666666
$ws.Range("H8").Value = 83.32245499999999
$ws.Range("I8").Value = 0.1759264889641896
$ws.Range("J8").Value = 0.1759264889641896
$ws.Range("M8").Value = 39.327127
$ws.Range("N8").Value = 117.981381
$ws.Range("O8").Value = 0.9923865713449503
$ws.Range("P8").Value = 0.9923865713449502
$ws.Range("Q8").Value = 1092.277589912261
$ws.Range("R8").Value = 9830.498309210354
$ws.Range("S8").Value = 0.1745870851919274
$ws.Range("T8").Value = 0.1745870851919273

$ws.Range("G9").Value = 27.77415166666666
$ws.Range("H9").Value = 83.32245499999999
$ws.Range("I9").Value = 0.1759264889641896
$ws.Range("J9").Value = 0.1759264889641896
$ws.Range("O9").Value = 0.001455135597170125
$ws.Range("P9").Value = 0.001455135597170125
$ws.Range("S9").Value = 0.0002559968965769495
$ws.Range("T9").Value = 0.0002559968965769494

$ws.Range("G10").Value = 27.77415166666666
$ws.Range("H10").Value = 83.32245499999999
$ws.Range("I10").Value = 0.1759264889641896
$ws.Range("J10").Value = 0.1759264889641896
$ws.Range("O10").Value = 0.00615829305787961
$ws.Range("P10").Value = 0.006158293057879609
$ws.Range("Q10").Value = 6.778170617643332
$ws.Range("R10").Value = 61.00353555878999
$ws.Range("S10").Value = 0.001083406875685303
$ws.Range("T10").Value = 0.001083406875685302
